# Update quantity (F) / value (G) cells for items whose stock quantity
# decreased, and the corresponding "Sub Total:"/"Grand Total:" (B) cells
# that roll those values up, per company group.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F27").Value2 = 16
$ws.Range("G27").Value2 = 3148.32
$ws.Range("F49").Value2 = 77
$ws.Range("G49").Value2 = 20076.21
$ws.Range("B54").Value2 = 84537.28999999999
$ws.Range("F94").Value2 = 1
$ws.Range("G94").Value2 = 70.31999999999999
$ws.Range("F96").Value2 = 26
$ws.Range("G96").Value2 = 1649.18
$ws.Range("F97").Value2 = 180
$ws.Range("G97").Value2 = 11466
$ws.Range("B116").Value2 = 155913.55
$ws.Range("F151").Value2 = 51
$ws.Range("G151").Value2 = 7315.44
$ws.Range("F163").Value2 = 352
$ws.Range("G163").Value2 = 34077.12
$ws.Range("F165").Value2 = 126
$ws.Range("G165").Value2 = 13173.3
$ws.Range("B166").Value2 = 113718.87
$ws.Range("F181").Value2 = 12
$ws.Range("G181").Value2 = 1448.28
$ws.Range("F185").Value2 = 1
$ws.Range("G185").Value2 = 72.34999999999999
$ws.Range("B188").Value2 = 5959.59
$ws.Range("F217").Value2 = 12
$ws.Range("G217").Value2 = 84.48
$ws.Range("B221").Value2 = 44862.11
$ws.Range("F223").Value2 = 130
$ws.Range("G223").Value2 = 14859
$ws.Range("B229").Value2 = 60185.14
$ws.Range("F322").Value2 = 185
$ws.Range("G322").Value2 = 8813.4
$ws.Range("F345").Value2 = 30
$ws.Range("G345").Value2 = 14831.4
$ws.Range("B360").Value2 = 164903.29
$ws.Range("F364").Value2 = 1
$ws.Range("G364").Value2 = 446.65
$ws.Range("B368").Value2 = 15010.75
$ws.Range("F395").Value2 = 66
$ws.Range("G395").Value2 = 2438.7
$ws.Range("F398").Value2 = 350
$ws.Range("G398").Value2 = 49206.5
$ws.Range("B400").Value2 = 63836.57
$ws.Range("F438").Value2 = 43
$ws.Range("G438").Value2 = 4224.75
$ws.Range("F447").Value2 = 5
$ws.Range("G447").Value2 = 932.2
$ws.Range("B456").Value2 = 100828.52
$ws.Range("F465").Value2 = 329
$ws.Range("G465").Value2 = 4326.35
$ws.Range("F466").Value2 = 449
$ws.Range("G466").Value2 = 5751.69
$ws.Range("F471").Value2 = 307
$ws.Range("G471").Value2 = 6057.11
$ws.Range("B482").Value2 = 73365.07000000001
$ws.Range("F491").Value2 = 306
$ws.Range("G491").Value2 = 18864.9
$ws.Range("F493").Value2 = 193
$ws.Range("G493").Value2 = 8932.040000000001
$ws.Range("B504").Value2 = 125826.78
$ws.Range("F522").Value2 = 19
$ws.Range("G522").Value2 = 2796.8
$ws.Range("F523").Value2 = 38
$ws.Range("G523").Value2 = 6125.6
$ws.Range("F524").Value2 = 59
$ws.Range("G524").Value2 = 12750.49
$ws.Range("B526").Value2 = 30053.41
$ws.Range("F557").Value2 = 94
$ws.Range("G557").Value2 = 1966.48
$ws.Range("B573").Value2 = 55713.9
$ws.Range("F601").Value2 = 58
$ws.Range("G601").Value2 = 9307.84
$ws.Range("F611").Value2 = 78
$ws.Range("G611").Value2 = 2137.2
$ws.Range("F616").Value2 = 84
$ws.Range("G616").Value2 = 5333.16
$ws.Range("F617").Value2 = 92
$ws.Range("G617").Value2 = 6424.36
$ws.Range("F618").Value2 = 10
$ws.Range("G618").Value2 = 1417.5
$ws.Range("B623").Value2 = 118625.28
$ws.Range("F644").Value2 = 14
$ws.Range("G644").Value2 = 604.52
$ws.Range("B646").Value2 = 11530.73
$ws.Range("F672").Value2 = 81
$ws.Range("G672").Value2 = 6714.09
$ws.Range("F673").Value2 = 77
$ws.Range("G673").Value2 = 9382.450000000001
$ws.Range("B680").Value2 = 106314.41
$ws.Range("F695").Value2 = 21
$ws.Range("G695").Value2 = 1912.68
$ws.Range("B696").Value2 = 1912.68
$ws.Range("F704").Value2 = 9
$ws.Range("G704").Value2 = 734.04
$ws.Range("F716").Value2 = 88
$ws.Range("G716").Value2 = 10622.48
$ws.Range("B718").Value2 = 48228.23
$ws.Range("F722").Value2 = 17
$ws.Range("G722").Value2 = 2786.13
$ws.Range("F730").Value2 = 17
$ws.Range("G730").Value2 = 897.09
$ws.Range("F732").Value2 = 32
$ws.Range("G732").Value2 = 5029.44
$ws.Range("B748").Value2 = 52137.25
$ws.Range("F786").Value2 = 1952
$ws.Range("G786").Value2 = 318390.72
$ws.Range("F787").Value2 = 190
$ws.Range("G787").Value2 = 53745.3
$ws.Range("F789").Value2 = 40
$ws.Range("G789").Value2 = 1525.6
$ws.Range("B794").Value2 = 446923.24
$ws.Range("B799").Value2 = 3085595.12
$ws.Range("B800").Value2 = 3085595.12
